$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.188000000000001
$ws.Range("A3").Value = -21.453
$ws.Range("A14").Value = -21.5
$ws.Range("A16").Value = -21.736
$ws.Range("B18").Value = 5.494
$ws.Range("A21").Value = -20.326
$ws.Range("A23").Value = -20.673
$ws.Range("B24").Value = 6.055
$ws.Range("A25").Value = -21.521
$ws.Range("B25").Value = 6.936
$ws.Range("A26").Value = -21.481
$ws.Range("B27").Value = 6.529999999999999
$ws.Range("A29").Value = -21.104
$ws.Range("B30").Value = 6.738000000000001
$ws.Range("B31").Value = 6.545
$ws.Range("B39").Value = 7.582000000000001
$ws.Range("A40").Value = -20.404
$ws.Range("B42").Value = 8.33
$ws.Range("B48").Value = 5.274
$ws.Range("B51").Value = 6.802
$ws.Range("B52").Value = 5.813000000000001
$ws.Range("A53").Value = -21.924
$ws.Range("B55").Value = 4.670999999999999
$ws.Range("B56").Value = 4.987
$ws.Range("A57").Value = -21.607
$ws.Range("B57").Value = 6.21
$ws.Range("A59").Value = -22.404
$ws.Range("B60").Value = 5.896999999999999
$ws.Range("A65").Value = -21.418
$ws.Range("A69").Value = -21.507
$ws.Range("B73").Value = 6.804
$ws.Range("B74").Value = 8.995000000000001
$ws.Range("A79").Value = -21.246
$ws.Range("A83").Value = -21.351
$ws.Range("B89").Value = 5.787999999999999
$ws.Range("B90").Value = 5.767
$ws.Range("A91").Value = -21.533
$ws.Range("B92").Value = 5.898999999999999
$ws.Range("A93").Value = -21.324
$ws.Range("A100").Value = -21.587
